# LMS-2391 Split TSV files where necessary. Switched tab to spaces in all python code.
# Re-creates the OD600-Example.xlsx edit:
#  - append a duplicate data row (row 6) to the "openbis-data" sheet, copied from row 2
#  - make "openbis-data" the active/selected sheet (was "openbis-metadata")
#  - select the new row (A6:XFD6) on "openbis-data"
#  - add a pageSetup (paper size / orientation) to "openbis-data"

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item(1)   # openbis-metadata
$wsData = $wb.Worksheets.Item(2)   # openbis-data

# ---- Copy row 2 into new row 6 on the openbis-data sheet ----
$srcRow = 2
$dstRow = 6
$lastCol = 21   # column U

for ($c = 1; $c -le $lastCol; $c++) {
    $srcCell = $wsData.Cells.Item($srcRow, $c)
    $dstCell = $wsData.Cells.Item($dstRow, $c)
    $dstCell.Value = $srcCell.Value()
}

# ---- Make openbis-data the active sheet and select the new row ----
$wsData.Activate() | Out-Null
$wsData.Rows("6:6").Select() | Out-Null

# ---- Add page setup info (paper size + orientation) to openbis-data ----
$ps = $wsData.PageSetup
$ps.PaperSize = 10
$ps.Orientation = 1
